$d = $word.ActiveDocument

function Set-ParagraphRed($searchText) {
    $found = $d.Content.Duplicate
    $ok = $found.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        return
    }
    $targetStart = $found.Start

    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $targetStart) {
            $p.Range.Font.Color = 255
            break
        }
    }
}

Set-ParagraphRed("Write Unit Tests for your logic, controllers, actions, helpers, etc.")
Set-ParagraphRed("You should cover at least 80% of your business logic.")
